$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column A, shifting everything right.
$ws.Columns.Item(1).Insert()

# New column A: "Match ID" header in row 3, value 24 for data rows 4-20,
# and 24 (no special style) in the hidden total row 21.
$ws.Cells.Item(3, 1).Value = "Match ID"

for ($r = 4; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = 24
}
$ws.Cells.Item(21, 1).Value = 24
$ws.Rows.Item(21).AutoFit()

# Header + data cells in the new column get a bold font (no border).
$ws.Range("A3:A20").Font.Bold = $true

# Update selection to match the new focus on the Match ID column.
$ws.Range("A3:A20").Select()
